# Retaking experimental results to deal with phase-shift error:
# append " - 50" (trial marker) to the Vds measurement labels in row 2-5,
# columns B and C of the dataset sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Vds Bot (Pos) - H1 - 50"
$ws.Range("C2").Value = "Vds Bot (Pos) - H1 - 50"

$ws.Range("B3").Value = "Vds Bot (Neg) - H2 -50"
$ws.Range("C3").Value = "Vds Bot (Neg) - H2 -50"

$ws.Range("B4").Value = "Vds Top(Neg) - H1 - 50"
$ws.Range("C4").Value = "Vds Bot (Pos) - H1 - 50"

$ws.Range("B5").Value = "VdsTop (Pos) - H2 - 50"
$ws.Range("C5").Value = "Vds Bot (Neg) - H2 - 50"

# Update the active selection on the sheet to C6, matching the saved view state.
$ws.Range("C6").Select()
